# Apply DB-diagram re-order: APUESTA -> EVENTO -> MERCADO
#                          becomes APUESTA -> MERCADO -> EVENTO
#
# Concretely, in the APUESTA block:
#   - field "Id_Evento: int" becomes "Id_Mercado: int"
#   - "Clave ajena: Id_Evento (EVENTO)" becomes "Clave ajena: Id_Mercado (MERCADO)"
# (the MERCADO block's own, unrelated "Id_Evento" foreign key to EVENTO is left
# untouched) and the hidden "_GoBack" bookmark (last cursor position) moves
# from just before the EVENTO block (which used to immediately follow APUESTA)
# to the empty paragraph right after the MERCADO block's "Clave ajena:
# Id_Evento (EVENTO)" line (EVENTO is now the last block in the document).

$d = $word.ActiveDocument
$count = $d.Paragraphs.Count

# --- find the section-header paragraphs (APUESTA, EVENTO, MERCADO) ---------
$apuestaHeader = 0
$eventoHeader = 0
$mercadoHeader = 0

for ($i = 1; $i -le $count; $i++) {
    $text = $d.Paragraphs.Item($i).Range.Text.Trim()
    if ($text -eq "APUESTA") { $apuestaHeader = $i }
    elseif ($text -eq "EVENTO") { $eventoHeader = $i }
    elseif ($text -eq "MERCADO") { $mercadoHeader = $i }
}

# --- within the APUESTA block (between the APUESTA and EVENTO headers), ----
# --- find the field-list line and the "Clave ajena" line -------------------
$apuestaFieldPara = $null
$apuestaFkPara = $null

for ($i = $apuestaHeader; $i -lt $eventoHeader; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text.Trim()

    if ($text -eq "Id_Evento: int") {
        $apuestaFieldPara = $para
    }
    elseif ($text -eq "Clave ajena: Id_Evento (EVENTO)") {
        $apuestaFkPara = $para
    }
}

# --- 1) "Id_Evento" -> "Id_Mercado" in the APUESTA field list --------------
$apuestaFieldPara.Range.Find.Execute("Id_Evento", $true, $false, $false, $false, $false,
                                      $true, 1, $false, "Id_Mercado", 2)

# --- 2) "Clave ajena: Id_Evento (EVENTO)" -> "Clave ajena: Id_Mercado (MERCADO)" --
#     rebuilt as "Id_Mercado" + " (MERCADO" + ")" (three separate runs, matching
#     the target markup) instead of a single in-place text swap.
$apuestaFkPara.Range.Find.Execute("Id_Evento", $true, $false, $false, $false, $false,
                                   $true, 1, $false, "Id_Mercado", 2)
$apuestaFkPara.Range.Find.Execute(" (EVENTO)", $true, $false, $false, $false, $false,
                                   $true, 1, $false, " (MERCADO", 2)
$apuestaFkPara.Range.InsertAfter(")")

# --- 3)/4) move the "_GoBack" bookmark -------------------------------------
# Within the MERCADO block (from the MERCADO header to the end of the
# document), find its own "Clave ajena: Id_Evento (EVENTO)" line; the bookmark
# goes on the empty paragraph right after it.
$mercadoFkPara = $null
for ($i = $mercadoHeader; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.Trim() -eq "Clave ajena: Id_Evento (EVENTO)") {
        $mercadoFkPara = $para
    }
}

$target = $mercadoFkPara.Next()
$d.Bookmarks.Add("_GoBack", $target.Range)
